$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (rows 2-15): break from the SIN shared-formula chain and
# instead compute COS() of the value to its immediate left (column C).
$ws.Range("D2").Formula = "=COS(C2)"
$ws.Range("D3:D15").Formula = "=COS(C3)"

# Column G (rows 2-15): same idea, COS() of column F.
$ws.Range("G2").Formula = "=COS(F2)"
$ws.Range("G3:G15").Formula = "=COS(F3)"

# Append two new rows (16 and 17) of literal 1s across columns A:G.
$ws.Range("A16:G16").Value = 1
$ws.Range("A17:G17").Value = 1

# Update the selected cell to reflect where the user ended up.
$ws.Range("H17").Select()

# Page setup was touched as well (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
